function Delete-ParaMarks($doc, $pos, $count) {
    # Paragraph-mark-only ranges only ever collapse one mark per Delete()
    # call in this environment, so remove them one at a time.
    for ($i = 0; $i -lt $count; $i++) {
        $r = $doc.Range($pos, $pos + 1)
        $r.Delete()
    }
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: paragraph "[  ]  -  Create Report ... - Report 1  CURRENTLY
# WORKING ON" / empty paragraph / "[  ]  -  Generate JavaDocs- will do this
# after reports are written" collapse into a single paragraph with new text.
# ---------------------------------------------------------------------------

# 1a. Remove the " " + Wingdings-symbol + " CURRENTLY WORKING ON" text that
#     trails " - Report 1" (keep the paragraph mark itself for now).
$full = $d.Content.Text
$anchor = " - Report 1"
$startDelete = $full.IndexOf($anchor) + $anchor.Length
$paraMarkPos = $d.Paragraphs.Item(3).Range.End - 1
$r = $d.Range($startDelete, $paraMarkPos)
$r.Text = ""

# 1b. Delete the two paragraph marks that separate the "Report 1" paragraph,
#     the blank paragraph after it, and the "Generate JavaDocs" paragraph, so
#     all three collapse into one paragraph.
$full = $d.Content.Text
$anchor = " - Report 1"
$pos = $full.IndexOf($anchor) + $anchor.Length
Delete-ParaMarks $d $pos 2

# 1c. Remove the now-adjacent "[  ]  -  Generate " lead-in text, leaving the
#     "JavaDocs"/proofErr run untouched immediately after " - Report 1".
$full = $d.Content.Text
$anchor = " - Report 1"
$startDelete = $full.IndexOf($anchor) + $anchor.Length
$endDelete = $full.IndexOf("JavaDocs")
$r = $d.Range($startDelete, $endDelete)
$r.Text = ""

# 1d. Swap "JavaDocs" for "ComboBox" (the spellStart/spellEnd wrapper stays
#     put since only the run's text content changes).
$d.Content.Find.Execute("JavaDocs", $true, $false, $false, $false, $false, `
    $true, 1, $false, "ComboBox", 2) | Out-Null

# 1e. Replace the old trailing remark with the new continuation text.
$d.Content.Find.Execute("- will do this after reports are written", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    " choice as you suggested, then will write the code to pull these from the database and generate the report. After I’m done with this I should be ready for submission.", `
    2) | Out-Null

# 1f. Insert the new sentence right after " - Report 1".
$full = $d.Content.Text
$anchor = " - Report 1"
$pos = $full.IndexOf($anchor) + $anchor.Length
$r = $d.Range($pos, $pos)
$r.InsertAfter(". I am changing my controller to the ")

# ---------------------------------------------------------------------------
# Change 2: delete the "[X]  -  View Schedule by Contact ID - Report 2  Need
# to figure out what to pass into initialize method" paragraph along with one
# of the blank paragraphs around it.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$anchor = "Completed But Not Working Properly"
$afterProperly = $full.IndexOf($anchor) + $anchor.Length

# 2a. Delete the blank paragraph mark immediately after "Not Working Properly".
Delete-ParaMarks $d $afterProperly 1

# 2b. Clear the "[X]  -  View Schedule ..." paragraph's text.
$full = $d.Content.Text
$start = $full.IndexOf("[X]")
$end = $start + "[X]  -  View Schedule by Contact ID - Report 2  Need to figure out what to pass into initialize method".Length
$r = $d.Range($start, $end)
$r.Text = ""

# 2c. Delete the now-empty paragraph's mark so it collapses away completely.
$full = $d.Content.Text
$start = $full.IndexOf("[X]")
if ($start -lt 0) { $start = $afterProperly }
Delete-ParaMarks $d $start 1
